# Automatische test-sync: 2025-06-19 22:07:50
#
# 1. Append a new incoming mail-log entry (row 36) to the "Logs" sheet.
# 2. Re-order the "Klacht/Probleem", "IT/Technisch probleem" and
#    "Factuur/Administratie" rows on the "Dashboard" sheet and append a
#    new "Overig" category row (row 11).
# 3. Grow the conditional-formatting ranges on "Logs" to cover the new row.
# 4. Grow the bar-chart's category/value series references on "Dashboard"
#    to include the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append new row 36
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(36, 1).Value = "Uitnodiging voor netwerkevent"
$logs.Cells.Item(36, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(36, 3).Value = "Graag nodig ik u uit voor ons zakelijke netwerkevent volgende maand."
$logs.Cells.Item(36, 4).Value = "Overig"
$logs.Cells.Item(36, 6).Value = "2025-06-19 22:07:35"
$logs.Cells.Item(36, 7).Value = "Nee"

# ---------------------------------------------------------------------
# 2. Dashboard sheet: re-order categories + append "Overig" row
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(5, 1).Value = "Klacht / Probleem"
$dash.Cells.Item(6, 1).Value = "IT / Technisch probleem"
$dash.Cells.Item(7, 1).Value = "Factuur / Administratie"

$dash.Cells.Item(11, 1).Value = "Overig"
$dash.Cells.Item(11, 2).Value = 1

# ---------------------------------------------------------------------
# 3. Logs sheet: extend conditional formatting ranges to row 36
# ---------------------------------------------------------------------
$dRule = $logs.Range("D2:D35").FormatConditions.Item(1)
$dRule.ModifyAppliesToRange($logs.Range("D2:D36"))

$gRule = $logs.Range("G2:G35").FormatConditions.Item(1)
$gRule.ModifyAppliesToRange($logs.Range("G2:G36"))

# ---------------------------------------------------------------------
# 4. Dashboard chart: extend category/value series to row 11
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$11,'Dashboard'!`$B`$2:`$B`$11,1)"
